# Updated symbol list (refreshed crypto price/volume snapshot).
# Column D ("Price") and some Column E ("Volume(1h)") cells are stored as
# plain text (t="inlineStr"/shared-string), not numbers, in this sheet.
# Assigning a numeric-looking string straight to .Value makes Excel coerce
# it to a real number (and drop formatting like trailing zeros), so for
# column D we prefix the new value with a leading apostrophe to force a
# text literal, then reset .Style back to "Normal" so the quote-prefix
# style Excel applies doesn't leave a stray style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'249.00"
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.Value = "'5.364"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'0.05613"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'3.408"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.Value = "'6.374"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.Value = "'0.9526"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.Value = "'0.07535"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.Value = "'0.03200"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = "'0.03092"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'0.09307"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.Value = "'3.552"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = "'0.001602"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = "'0.04706"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.Value = "'0.0005769"
$cell.Style = "Normal"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$cell = $ws.Range("D19")
$cell.Value = "'0.006280"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = "'0.0001498"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.Value = "'3.754"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = "'2.148"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.Value = "'0.3300"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.Value = "'0.1313"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = "'0.03955"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = "'0.007016"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = "'0.1066"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'0.003112"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.Value = "'0.008761"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'0.00005728"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.Value = "'0.00000000749"
$cell.Style = "Normal"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
$cell = $ws.Range("D48")
$cell.Value = "'0.6800"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.Value = "'0.1708"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'0.00002097"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = "'0.01009"
$cell.Style = "Normal"
